$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$newValues = @(
    "94-25=69",
    "53-9=44",
    "39+54=93",
    "74+17=91",
    "84-68=16",
    "54-39=15",
    "91-12=79",
    "80-31=49",
    "91-54=37",
    "80-43=37",
    "46-37=9",
    "82-69=13",
    "44-8=36",
    "4+19=23",
    "29+42=71",
    "90-13=77",
    "92-45=47",
    "9+78=87",
    "9+24=33",
    "36+27=63",
    "86-7=79",
    "63-57=6",
    "37+39=76",
    "50-35=15",
    "60-55=5",
    "84-28=56",
    "61-6=55",
    "80-2=78",
    "75-36=39",
    "65-37=28",
    "31-6=25",
    "63-6=57",
    "62-13=49",
    "16+56=72",
    "49+32=81",
    "70-48=22",
    "16+25=41",
    "25-17=8",
    "30-9=21",
    "17+6=23",
    "28+26=54",
    "43-6=37",
    "55+38=93",
    "70-58=12",
    "34+29=63",
    "95-67=28",
    "61-44=17",
    "70-43=27",
    "89+4=93",
    "60-31=29",
    "96-18=78",
    "87+9=96",
    "31-13=18",
    "79+15=94",
    "9+55=64",
    "90-7=83",
    "78+8=86",
    "63-19=44",
    "82-79=3",
    "37+25=62",
    "58+36=94",
    "92-15=77",
    "43+19=62",
    "37+45=82",
    "70-62=8",
    "42-4=38",
    "46-17=29",
    "3+59=62",
    "12+19=31",
    "29+65=94",
    "7+44=51",
    "94-6=88",
    "74-49=25",
    "59+38=97",
    "48+19=67",
    "46-8=38",
    "80-7=73",
    "29+54=83",
    "11-2=9",
    "43-7=36",
    "88+4=92",
    "55+6=61",
    "66-19=47",
    "41-22=19",
    "48+24=72",
    "61-5=56",
    "27+26=53",
    "22-18=4",
    "38+39=77",
    "9+52=61",
    "73-6=67",
    "39+47=86",
    "59+27=86",
    "13-5=8",
    "59+19=78",
    "9+44=53",
    "58+29=87",
    "19+13=32",
    "8+14=22",
    "34+28=62"
)

$rows = $t.Rows.Count
$cols = $t.Columns.Count
$idx = 0
for ($r = 1; $r -le $rows; $r++) {
  for ($c = 1; $c -le $cols; $c++) {
    $cell = $t.Cell($r, $c)
    $cell.Range.Text = $newValues[$idx]
    $idx = $idx + 1
  }
}

Write-Host "Updated $idx cells"
